# "fixed export and fixing maps"
#
# The table is being trimmed down from a 3-year (1989 / 2002 / 2014) export
# to a single 2014-only export, and the extra descriptive caption row is
# dropped:
#   - Row 2 ("(according to the population census data)") is removed
#     entirely, shifting everything below it up by one row.
#   - Columns B:C (the 1989 and 2002 data columns) are removed entirely,
#     leaving only the 2014 column, which shifts left into column B.
#   - Four blank rows are appended at the bottom (new rows 6-9).
#   - All rows get a uniform custom row height of 20.1 points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "(according to the population census data)" caption row.
$ws.Rows("2").Delete()

# Keep only the 2014 column; drop the 1989/2002 columns (old B:C).
$ws.Range("B1:C1").EntireColumn.Delete()

# Add four trailing blank rows, fully empty (no leftover values/formats).
$ws.Rows("6:9").Insert()
$ws.Rows("6:9").Clear()

# Uniform custom row height across the whole table.
$ws.Rows("1:9").RowHeight = 20.1
